# Append: 2026-02-04 13:08 JST
# Update the "取得日時" (retrieved at) timestamp in column A for every
# existing data row on the first sheet ("ランサーズ") from
# "2026-02-04 12:54:41" to "2026-02-04 13:08:28".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldTimestamp = "2026-02-04 12:54:41"
$newTimestamp = "2026-02-04 13:08:28"

# Find the last used row based on column A, then update every row (from
# row 2, skipping the header row) whose timestamp matches the old value.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
